# Frontend Functionalities.xlsx — add two new functionality rows:
#   22 -> "User Preview"          (red fill, like row 10 "Add product")
#   23 -> "User update username"  (green fill, like row 2 "User Login")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 23: index 22, "User Preview"
$ws.Range("A23").Value = 22
$ws.Range("B23").Value = "User Preview"

# New row 24: index 23, "User update username"
$ws.Range("A24").Value = 23
$ws.Range("B24").Value = "User update username"

# Match the C-column cell formatting (fill + border) used elsewhere in the
# table: red status style (same as row 10) for the new "User Preview" row,
# green status style (same as row 2) for "User update username".
$ws.Range("C10").Copy()
$ws.Range("C23").PasteSpecial(-4122)

$ws.Range("C2").Copy()
$ws.Range("C24").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Move the active selection, matching the saved view state.
$ws.Range("G12").Select()
